# Insert a new column "search_term_selic" before the existing
# "search_term_igpm" column (column D), shifting the igpm header + its
# data one column to the right (to column E), leaving the new column D
# blank apart from its header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift column D (and everything to its right) one column to the right.
# This moves the old D1 header ("search_term_igpm") and its 239 data
# values down into column E, and leaves a blank column D behind
# (inheriting the column's existing number formatting/style).
$ws.Columns.Item(4).Insert()

# Give the newly-inserted, now-blank column D its new header text.
$ws.Cells.Item(1, 4).Value = "search_term_selic"
